$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C9 was the text "4h40 for the moment" - replace it with the numeric value 6
$ws.Range("C9").Value = 6

# Update the active selection shown in the sheet view to D12
$ws.Range("D12").Select()
